# Weekly update: insert two new price-report rows (Terminal Hortofrutícola
# Agro Chillán - Cilantro) ahead of the existing row 62, shifting the rest
# of the table down by two rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows at 62-63 (existing rows 62..77 shift down to 64..79)
$ws.Range("A62:A63").EntireRow.Insert()

# New row 62: Cilantro, Primera
$ws.Cells.Item(62, 1).Value = 7
$ws.Cells.Item(62, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(62, 3).Value = "Ñuble"
$ws.Cells.Item(62, 4).Value = 44798
$ws.Cells.Item(62, 5).Value = 16
$ws.Cells.Item(62, 6).Value = 100112040
$ws.Cells.Item(62, 7).Value = "Cilantro"
$ws.Cells.Item(62, 8).Value = "Sin especificar"
$ws.Cells.Item(62, 9).Value = "Primera"
$ws.Cells.Item(62, 10).Value = 200
$ws.Cells.Item(62, 11).Value = 700
$ws.Cells.Item(62, 12).Value = 800
$ws.Cells.Item(62, 13).Value = 750
$ws.Cells.Item(62, 14).Value = "$/atado 0,5 a 1 kilo"
$ws.Cells.Item(62, 15).Value = "Provincia de Diguillín"
$ws.Cells.Item(62, 16).Value = 750
$ws.Cells.Item(62, 17).Value = 1
$ws.Cells.Item(62, 18).Value = "Hortaliza"

# New row 63: Cilantro, Segunda
$ws.Cells.Item(63, 1).Value = 7
$ws.Cells.Item(63, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(63, 3).Value = "Ñuble"
$ws.Cells.Item(63, 4).Value = 44798
$ws.Cells.Item(63, 5).Value = 16
$ws.Cells.Item(63, 6).Value = 100112040
$ws.Cells.Item(63, 7).Value = "Cilantro"
$ws.Cells.Item(63, 8).Value = "Sin especificar"
$ws.Cells.Item(63, 9).Value = "Segunda"
$ws.Cells.Item(63, 10).Value = 150
$ws.Cells.Item(63, 11).Value = 600
$ws.Cells.Item(63, 12).Value = 600
$ws.Cells.Item(63, 13).Value = 600
$ws.Cells.Item(63, 14).Value = "$/atado 0,5 a 1 kilo"
$ws.Cells.Item(63, 15).Value = "Provincia de Diguillín"
$ws.Cells.Item(63, 16).Value = 600
$ws.Cells.Item(63, 17).Value = 1
$ws.Cells.Item(63, 18).Value = "Hortaliza"
